$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$row2 = @{
    "G2" = 3.1
    "H2" = 3.6
    "I2" = 2.15
    "J2" = 3.75
    "K2" = 2.2
    "L2" = 2.75
    "M2" = 1.04
    "N2" = 12
    "O2" = 1.25
    "P2" = 3.75
    "Q2" = 1.83
    "R2" = 1.98
    "S2" = 1.36
    "T2" = 3
    "U2" = 1.67
    "V2" = 2.1
    "W2" = 11
    "X2" = 17
    "Z2" = 34
    "AA2" = 23
    "AB2" = 29
    "AC2" = 12
    "AD2" = 7
    "AE2" = 13
    "AF2" = 41
    "AG2" = 201
    "AH2" = 8.5
    "AI2" = 11
    "AJ2" = 9
    "AK2" = 19
    "AL2" = 17
    "AM2" = 26
    "AN2" = 5.5
    "AO2" = 17
    "AP2" = 23
    "AQ2" = 51
    "AR2" = 67
    "AS2" = 151
    "AT2" = 3
    "AU2" = 7.5
    "AV2" = 51
    "AW2" = 4.33
    "AX2" = 12
    "AY2" = 21
    "AZ2" = 41
    "BA2" = 51
    "BB2" = 126
}
foreach ($key in $row2.Keys) {
    $ws.Range($key).Value = $row2[$key]
}

# Row 3
$row3 = @{
    "G3" = 1.44
    "H3" = 4.35
    "I3" = 6.1
    "K3" = 2.37
    "L3" = 5.7
    "N3" = 13.3
    "O3" = 1.12
    "P3" = 4.45
    "Q3" = 1.55
    "R3" = 2.15
    "S3" = 1.28
    "T3" = 3.46
    "U3" = 1.7
    "V3" = 1.91
    "W3" = 8.25
    "X3" = 7.6
    "Y3" = 8.25
    "Z3" = 10.25
    "AA3" = 10.75
    "AB3" = 22
    "AC3" = 14
    "AD3" = 8.75
    "AE3" = 17
    "AF3" = 70
    "AG3" = 450
    "AH3" = 18.5
    "AJ3" = 19.5
    "AK3" = 120
    "AM3" = 55
    "AN3" = 3.35
    "AO3" = 6.5
    "AP3" = 14.5
    "AQ3" = 18
    "AR3" = 40
    "AS3" = 175
    "AT3" = 3.15
    "AU3" = 7.6
    "AV3" = 65
    "AW3" = 7.6
    "AX3" = 35
    "AY3" = 35
    "AZ3" = 200
    "BA3" = 250
    "BB3" = 400
}
foreach ($key in $row3.Keys) {
    $ws.Range($key).Value = $row3[$key]
}

# Row 4
$row4 = @{
    "G4" = 2.15
    "H4" = 3.55
    "I4" = 2.92
    "J4" = 2.67
    "K4" = 2.25
    "L4" = 3.35
    "P4" = 3.96
    "Q4" = 1.55
    "R4" = 2.15
    "S4" = 1.3
    "T4" = 3.32
    "U4" = 1.5
    "V4" = 2.27
    "W4" = 10.5
    "X4" = 12.5
    "Y4" = 8.75
    "Z4" = 22
    "AA4" = 15.5
    "AB4" = 20
    "AC4" = 14.5
    "AD4" = 7.3
    "AE4" = 11.5
    "AF4" = 40
    "AG4" = 250
    "AH4" = 12.5
    "AI4" = 18
    "AJ4" = 10.5
    "AK4" = 37
    "AL4" = 22
    "AM4" = 24
    "AN4" = 4.3
    "AO4" = 10.75
    "AP4" = 16
    "AQ4" = 40
    "AR4" = 60
    "AS4" = 175
    "AT4" = 3.15
    "AU4" = 6.4
    "AV4" = 45
    "AW4" = 5.1
    "AX4" = 15
    "AY4" = 19
    "AZ4" = 65
    "BA4" = 80
    "BB4" = 200
}
foreach ($key in $row4.Keys) {
    $ws.Range($key).Value = $row4[$key]
}

# Row 7
$row7 = @{
    "G7" = 2.8
    "H7" = 3.2
    "I7" = 2.55
    "K7" = 1.91
    "L7" = 3.5
    "S7" = 1.57
    "T7" = 2.25
    "W7" = 6.5
    "AC7" = 7
    "AJ7" = 11
    "AL7" = 26
    "AP7" = 34
    "AT7" = 2.25
    "AX7" = 17
    "BB7" = 301
}
foreach ($key in $row7.Keys) {
    $ws.Range($key).Value = $row7[$key]
}

# Row 10
$row10 = @{
    "G10" = 1.65
    "H10" = 3.5
    "I10" = 4.5
    "J10" = 2.25
    "K10" = 2.3
    "L10" = 4.75
    "O10" = 1.2
    "P10" = 4.33
    "R10" = 2.1
    "S10" = 1.33
    "T10" = 3.25
    "U10" = 1.67
    "W10" = 8.5
    "X10" = 9
    "Y10" = 8.5
    "Z10" = 13
    "AA10" = 13
    "AB10" = 23
    "AC10" = 13
    "AD10" = 7
    "AE10" = 13
    "AF10" = 41
    "AG10" = 151
    "AH10" = 15
    "AI10" = 26
    "AJ10" = 15
    "AK10" = 51
    "AL10" = 34
    "AM10" = 41
    "AN10" = 3.75
    "AO10" = 8.5
    "AP10" = 19
    "AQ10" = 26
    "AR10" = 41
    "AS10" = 126
    "AT10" = 3.25
    "AU10" = 8
    "AV10" = 51
    "AW10" = 6.5
    "AX10" = 23
    "AY10" = 29
    "AZ10" = 81
    "BA10" = 81
    "BB10" = 151
}
foreach ($key in $row10.Keys) {
    $ws.Range($key).Value = $row10[$key]
}

# Row 11
$row11 = @{
    "G11" = 1.65
    "H11" = 3.45
    "I11" = 5.1
    "J11" = 2.18
    "K11" = 2.15
    "L11" = 5.1
    "N11" = 6.9
    "O11" = 1.32
    "P11" = 2.85
    "Q11" = 1.93
    "R11" = 1.7
    "S11" = 1.42
    "T11" = 2.66
    "U11" = 1.88
    "V11" = 1.72
    "W11" = 5.9
    "X11" = 7.1
    "Y11" = 8.25
    "Z11" = 12.5
    "AA11" = 14
    "AB11" = 30
    "AC11" = 8.75
    "AD11" = 6.8
    "AE11" = 17
    "AG11" = 800
    "AH11" = 13
    "AI11" = 30
    "AJ11" = 16.5
    "AL11" = 55
    "AM11" = 55
    "AN11" = 3.45
    "AO11" = 7.8
    "AP11" = 17
    "AQ11" = 25
    "AR11" = 55
    "AS11" = 200
    "AT11" = 2.77
    "AU11" = 7.3
    "AV11" = 65
    "AW11" = 6.7
    "AX11" = 29
    "AY11" = 30
    "AZ11" = 175
    "BA11" = 175
    "BB11" = 400
}
foreach ($key in $row11.Keys) {
    $ws.Range($key).Value = $row11[$key]
}

# Row 12
$row12 = @{
    "I12" = 2.05
    "L12" = 2.6
    "AA12" = 35
    "AK12" = 19
    "AL12" = 17
    "AN12" = 5.4
    "AO12" = 19.5
    "AT12" = 2.57
    "AV12" = 55
    "AY12" = 17.5
    "AZ12" = 37
    "BA12" = 65
    "BB12" = 200
}
foreach ($key in $row12.Keys) {
    $ws.Range($key).Value = $row12[$key]
}

# Row 13
$row13 = @{
    "G13" = 4.8
    "H13" = 4.2
    "I13" = 1.57
    "J13" = 4.5
    "K13" = 2.52
    "L13" = 2.05
    "M13" = 1.02
    "N13" = 10
    "O13" = 1.13
    "P13" = 5.2
    "Q13" = 1.4
    "R13" = 2.7
    "S13" = 1.24
    "T13" = 3.7
    "U13" = 1.47
    "V13" = 2.5
    "W13" = 23
    "X13" = 37
    "Y13" = 15.5
    "AB13" = 30
    "AC13" = 10
    "AD13" = 9.25
    "AG13" = 175
    "AH13" = 11.25
    "AI13" = 10.25
    "AJ13" = 8.25
    "AK13" = 13.5
    "AL13" = 11
    "AM13" = 16.5
    "AN13" = 7.2
    "AT13" = 3.7
    "AW13" = 3.95
    "AX13" = 7.5
    "AY13" = 12.5
    "AZ13" = 20
    "BA13" = 35
    "BB13" = 110
}
foreach ($key in $row13.Keys) {
    $ws.Range($key).Value = $row13[$key]
}

# Cells to clear (set blank, matching the target state)
$clearCells = @("M4","N4","M12","N12")
foreach ($key in $clearCells) {
    $ws.Range($key).ClearContents()
}

Write-Host "Applied odds update for week of 2024-12-11"